# feat: add 2022-Q3 data
#
# - Insert a brand-new worksheet "2022-Q3" right after "总计" (so it becomes
#   the 2nd tab), pushing the previously-existing quarter sheets one slot
#   later: 2022-Q2 -> position 3, 2022-Q1 -> position 4, 2021-Q4 -> position 5.
# - The new "2022-Q3" sheet carries the same single fund (540002, 汇丰晋信龙腾混合)
#   that used to be reported for 2022-Q2, refreshed with the new quarter's
#   metrics.
# - The "总计" (summary) sheet gets a new leading row for 2022-Q3, and the
#   existing rows shift down by one / get re-indexed in column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet before the current "2022-Q2" tab.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row (identical layout/style to every other quarter sheet).
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Copy header formatting (bold + border, style used by every quarter sheet)
# from the existing "2022-Q2" sheet so the new tab matches it exactly.
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)

# Data row.
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "540002"
$q3Sheet.Range("C2").Value = "汇丰晋信龙腾混合"
$q3Sheet.Range("D2").Value = "4.72"
$q3Sheet.Range("E2").Value = "93.98"
$q3Sheet.Range("F2").Value = "6.25"
$q3Sheet.Range("G2").Value = "0.2950"
$q3Sheet.Range("H2").Value = 4

$q2Sheet.Range("A2").Copy()
$q3Sheet.Range("A2").PasteSpecial(-4122)
$q3Sheet.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: new top row for 2022-Q3, shift the
#    previously-existing rows down by one (content-wise) and renumber the
#    index column.
#    NOTE: reading a cell's `.Value` into an expression requires the
#    method-call form `.Value()` in this host - the bare property getter
#    does not resolve when used as an rvalue.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("B5").Value = $totalSheet.Range("B4").Value()
$totalSheet.Range("C5").Value = $totalSheet.Range("C4").Value()
$totalSheet.Range("D5").Value = $totalSheet.Range("D4").Value()

$totalSheet.Range("B4").Value = $totalSheet.Range("B3").Value()
$totalSheet.Range("C4").Value = $totalSheet.Range("C3").Value()
$totalSheet.Range("D4").Value = $totalSheet.Range("D3").Value()

$totalSheet.Range("B3").Value = $totalSheet.Range("B2").Value()
$totalSheet.Range("C3").Value = $totalSheet.Range("C2").Value()
$totalSheet.Range("D3").Value = $totalSheet.Range("D2").Value()

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.3

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
